$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, pushing existing rows 60-80 down to 61-81
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new weekly price record
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44609
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100106
$ws.Range("H60").Value = "Oleaginosos"
$ws.Range("I60").Value = 100106002
$ws.Range("J60").Value = "Palta"
$ws.Range("K60").Value = "Hass"
$ws.Range("L60").Value = "Tercera"
$ws.Range("M60").Value = 200
$ws.Range("N60").Value = 68000
$ws.Range("O60").Value = 70000
$ws.Range("P60").Value = 69000
$ws.Range("Q60").Value = "$/caja 25 kilos"
$ws.Range("R60").Value = "Región de Coquimbo"
$ws.Range("S60").Value = 2760
$ws.Range("T60").Value = 25
